$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# CRS sheet: add the 4 new review rows (LH_REVIEW_CRS_001..004)
# ---------------------------------------------------------------------
$crs = $wb.Worksheets.Item("CRS")

# Cells are written in a specific order so the shared-strings table ends
# up built in the same sequence as the authored workbook.
$crs.Range("A2").Value = "LH_REVIEW_CRS_001"
$crs.Range("B2").Value = "Ahmed Abuzaid"
$crs.Range("C2").Value = "0ca4136"
$crs.Range("D2").Value = "1_the sheet name is named sheet 1`n2_the feature column not useful in this case"
$crs.Range("F2").Value = "omar sherif"

$crs.Range("A3").Value = "LH_REVIEW_CRS_002"
$crs.Range("D3").Value = "for the CRS item LH_CRS_005 – System Constrain, the current phrasing `n“web-based system/PC based” is a bit unclear"
$crs.Range("E3").Value = "can you make it more clear "

$crs.Range("E2").Value = "1_name the sheet according naming `nconvetion ""LH_CRS""`n2_we can suffice with the feature name in the id"

$crs.Range("A4").Value = "LH_REVIEW_CRS_003"
$crs.Range("A5").Value = "LH_REVIEW_CRS_004"

$crs.Range("D4").Value = "no comment on "" LH_CRS_ID-CONSTRAINS_006"""
$crs.Range("D5").Value = "no comment on ""LH_CRS_ADMIN-CONSTRAINS_007"""

$crs.Range("E4").Value = "no action"
$crs.Range("E5").Value = "no action"

# Remaining columns for rows 3-5 (reuse already-registered shared strings)
$crs.Range("B3").Value = "Ahmed Abuzaid"
$crs.Range("C3").Value = "0ca4136"
$crs.Range("F3").Value = "omar sherif"

$crs.Range("B4").Value = "Ahmed Abuzaid"
$crs.Range("C4").Value = "0ca4136"
$crs.Range("F4").Value = "omar sherif"

$crs.Range("B5").Value = "Ahmed Abuzaid"
$crs.Range("C5").Value = "0ca4136"
$crs.Range("F5").Value = "omar sherif"

# Status / Reviewer verification columns
$crs.Range("G2").Value = "open"
$crs.Range("H2").Value = "open"
$crs.Range("G3").Value = "open"
$crs.Range("H3").Value = "open"
$crs.Range("G4").Value = "closed"
$crs.Range("H4").Value = "closed"
$crs.Range("G5").Value = "closed"
$crs.Range("H5").Value = "closed"

# Wrap the long comment/action cells and size the rows to fit them
$crs.Range("D2").WrapText = $true
$crs.Range("E2").WrapText = $true
$crs.Range("D3").WrapText = $true

$crs.Rows.Item(2).RowHeight = 45
$crs.Rows.Item(3).RowHeight = 45

# Print page for the CRS sheet is switched to portrait orientation
$crs.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# SIQ sheet: row 2 (the long review comment) now uses the max row
# height of 409.5 (was 409.6), and the view/selection resets to A2.
# ---------------------------------------------------------------------
$siq = $wb.Worksheets.Item("SIQ")
$siq.Rows.Item(2).RowHeight = 409.5

# ---------------------------------------------------------------------
# View state: CRS becomes the active / selected sheet, SIQ's selection
# moves to A2 and it is no longer the active tab.
# ---------------------------------------------------------------------
[void]$siq.Activate()
[void]$siq.Range("A2").Select()

[void]$crs.Activate()
[void]$crs.Range("I8").Select()

Write-Host "done"
